$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.53 = 39373.25 pesos`n✅ 39373.25 pesos = 9.46 = 961.06 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 104.891
$wsTasas.Range("O10").Value = 4129.9
$wsTasas.Range("N12").Value = 4162
$wsTasas.Range("O12").Value = 101.59
